$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Section" header (E1) to "NRC", and add a new "Enseignant" header in F1
$ws.Range("E1").Value = "NRC"
$ws.Range("F1").Value = "Enseignant"

# Replace the old JS/MG "Section" codes in column E with numeric NRC values,
# and push the teacher name into the new column F
$ws.Range("E2").Value = 12345
$ws.Range("F2").Value = "J. Soucy"

$ws.Range("E3").Value = 12345
$ws.Range("F3").Value = "J. Soucy"

$ws.Range("E4").Value = 12345
$ws.Range("F4").Value = "J. Soucy"

$ws.Range("E5").Value = 12345
$ws.Range("F5").Value = "J. Soucy"

$ws.Range("E6").Value = 12345
$ws.Range("F6").Value = "J. Soucy"

$ws.Range("E7").Value = 54321
$ws.Range("F7").Value = "M. Genest"

$ws.Range("E8").Value = 54321
$ws.Range("F8").Value = "M. Genest"

$ws.Range("E9").Value = 54321
$ws.Range("F9").Value = "M. Genest"

$ws.Range("E10").Value = 54321
$ws.Range("F10").Value = "M. Genest"

$ws.Range("E11").Value = 54321
$ws.Range("F11").Value = "M. Genest"

# Update the active selection to match the target state
$ws.Range("G16").Select()
